# TS001_testscript.xlsx update — added new keyword switch_to_iframe.
# Rewrites the keyword-driven test script body (rows 1-20) with the new
# checkbox/radio-button jqueryui.com test flow, and relocates the trailing
# block of styled-but-empty D-column placeholder cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: write the brand-new cell text FIRST, in the exact order the ---
# --- strings are first introduced, so the shared-string table ends up   ---
# --- ordered the same way the authored workbook is.                    ---
$ws.Range("D3").Value = "To verify if the below:`n1. All checkboxes able to check."
$ws.Range("B4").Value = "Open the browser"
$ws.Range("C4").Value = "The browser opens successfully"
$ws.Range("D6").Value = "https://jqueryui.com/"
$ws.Range("B7").Value = "Open the checkbox radio page"
$ws.Range("C7").Value = "The checkbox radio page opens successfully"
$ws.Range("A8").Value = "click"
$ws.Range("B8").Value = "Checkboxradio link"
$ws.Range("C8").Value = "checkbox_page_css"
$ws.Range("B13").Value = "1radio"
$ws.Range("C13").Value = "ny_radio_css"
$ws.Range("B14").Value = "2radio"
$ws.Range("B15").Value = "3radio"
$ws.Range("C14").Value = "paris_radio_css"
$ws.Range("C15").Value = "london_radio_css"
$ws.Range("A12").Value = "switch_to_iframe"
$ws.Range("B12").Value = "iframe"
$ws.Range("C12").Value = "checkbox_iframe_css"
$ws.Range("C10").Value = "checkbox_no_icons_link_css"
$ws.Range("B10").Value = "checkboxnoicons link"
$ws.Range("B17").Value = "1check"
$ws.Range("C17").Value = "2star_check_css"
$ws.Range("B18").Value = "2check"
$ws.Range("B19").Value = "3check"
$ws.Range("B20").Value = "4check"
$ws.Range("C18").Value = "3star_check_css"
$ws.Range("C19").Value = "4star_check_css"
$ws.Range("C20").Value = "5star_check_css"
$ws.Range("B9").Value = "Click on no icons link"
$ws.Range("C9").Value = "The no icons page opens successfully"
$ws.Range("B11").Value = "Click on all radio buttons anc check"
$ws.Range("C11").Value = "Should be able to click on each radio button and that button is selected"
$ws.Range("B16").Value = "Click on all check boxes anc check"
$ws.Range("C16").Value = "Should be able to click on all check boxes and all should be checked."

# --- Step 2: (re)write the cells whose text already existed in the       ---
# --- workbook (header row + the Keyword column, which just repeats      ---
# --- "step"/"click"/etc.) so every row/column lines up with the new      ---
# --- layout.                                                             ---
$ws.Range("A1").Value = "Keyword"
$ws.Range("B1").Value = "Input1"
$ws.Range("C1").Value = "Input2"
$ws.Range("D1").Value = "Input3"
$ws.Range("A2").Value = "tc_id"
$ws.Range("D2").Value = "TS001"
$ws.Range("A3").Value = "tc_desc"
$ws.Range("A4").Value = "step"
$ws.Range("A5").Value = "open_browser"
$ws.Range("B5").Value = "Chrome"
$ws.Range("D5").Value = "Chrome"
$ws.Range("A6").Value = "enter_url"
$ws.Range("A7").Value = "step"
$ws.Range("A9").Value = "step"
$ws.Range("A10").Value = "click"
$ws.Range("A11").Value = "step"
$ws.Range("A13").Value = "click"
$ws.Range("A14").Value = "click"
$ws.Range("A15").Value = "click"
$ws.Range("A16").Value = "step"
$ws.Range("A17").Value = "click"
$ws.Range("A18").Value = "click"
$ws.Range("A19").Value = "click"
$ws.Range("A20").Value = "click"

# --- Step 3: row 3 (tc_desc) now wraps to two lines instead of four. ---
$ws.Rows.Item(3).RowHeight = 28.8

# --- Step 4: drop the old trailing placeholder rows (styled-only, blank ---
# --- D cells) entirely -- they get recreated a few rows further down.   ---
$ws.Range("D19").Clear()
$ws.Range("D21").Clear()
$ws.Range("D24").Clear()
$ws.Range("D26").Clear()
$ws.Range("D28").Clear()
$ws.Range("D30").Clear()
$ws.Range("D32").Clear()
$ws.Range("D33").Clear()
$ws.Range("D34").Clear()
$ws.Range("D36").Clear()
$ws.Range("D38").Clear()
$ws.Range("D39").Clear()
$ws.Range("D42").Clear()

# --- Step 5: recreate that placeholder block, shifted down by 5 rows. ---
$ws.Range("D24").WrapText = $true
$ws.Range("D26").WrapText = $true
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").WrapText = $true
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").WrapText = $true
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").WrapText = $true
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").WrapText = $true
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").WrapText = $true
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").WrapText = $true
$ws.Range("D47").NumberFormat = "@"

# --- Step 6: the selection left behind after the edit. ---
$ws.Range("C17").Select()
